{"js": "const replacements = [\n  [\"2024-07-28 Sunday\", \"2024-07-29 Monday\"],\n  [\"98\u00d753=5194\", \"28\u00d717=476\"],\n  [\"28\u00d789=2492\", \"16\u00d728=448\"],\n  [\"11\u00d721=231\", \"98\u00d779=7742\"],\n  [\"67\u00d738=2546\", \"18\u00d777=1386\"],\n  [\"35\u00d786=3010\", \"29\u00d732=928\"],\n  [\"15\u00d789=1335\", \"87\u00d798=8526\"],\n  [\"96\u00d728=2688\", \"49\u00d753=2597\"],\n  [\"68\u00d725=1700\", \"56\u00d773=4088\"],\n  [\"19\u00d785=1615\", \"94\u00d712=1128\"],\n  [\"20\u00d726=520\", \"28\u00d784=2352\"],\n  [\"29\u00d736=1044\", \"60\u00d739=2340\"],\n  [\"21\u00d735=735\", \"41\u00d737=1517\"],\n  [\"86\u00d778=6708\", \"30\u00d751=1530\"],\n  [\"14\u00d797=1358\", \"63\u00d743=2709\"],\n  [\"98\u00d755=5390\", \"43\u00d784=3612\"],\n  [\"74\u00d729=2146\", \"27\u00d769=1863\"],\n  [\"59\u00d752=3068\", \"46\u00d778=3588\"],\n  [\"59\u00d713=767\", \"32\u00d797=3104\"],\n  [\"64\u00d718=1152\", \"23\u00d769=1587\"],\n  [\"53\u00d744=2332\", \"50\u00d745=2250\"],\n  [\"27\u00d783=2241\", \"56\u00d782=4592\"],\n  [\"28\u00d757=1596\", \"30\u00d740=1200\"],\n  [\"23\u00d725=575\", \"54\u00d723=1242\"],\n  [\"44\u00d778=3432\", \"34\u00d740=1360\"],\n  [\"45\u00d798=4410\", \"17\u00d771=1207\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"2024-07-28 Sunday\", \"2024-07-29 Monday\"),\n  @(\"98\u00d753=5194\", \"28\u00d717=476\"),\n  @(\"28\u00d789=2492\", \"16\u00d728=448\"),\n  @(\"11\u00d721=231\", \"98\u00d779=7742\"),\n  @(\"67\u00d738=2546\", \"18\u00d777=1386\"),\n  @(\"35\u00d786=3010\", \"29\u00d732=928\"),\n  @(\"15\u00d789=1335\", \"87\u00d798=8526\"),\n  @(\"96\u00d728=2688\", \"49\u00d753=2597\"),\n  @(\"68\u00d725=1700\", \"56\u00d773=4088\"),\n  @(\"19\u00d785=1615\", \"94\u00d712=1128\"),\n  @(\"20\u00d726=520\", \"28\u00d784=2352\"),\n  @(\"29\u00d736=1044\", \"60\u00d739=2340\"),\n  @(\"21\u00d735=735\", \"41\u00d737=1517\"),\n  @(\"86\u00d778=6708\", \"30\u00d751=1530\"),\n  @(\"14\u00d797=1358\", \"63\u00d743=2709\"),\n  @(\"98\u00d755=5390\", \"43\u00d784=3612\"),\n  @(\"74\u00d729=2146\", \"27\u00d769=1863\"),\n  @(\"59\u00d752=3068\", \"46\u00d778=3588\"),\n  @(\"59\u00d713=767\", \"32\u00d797=3104\"),\n  @(\"64\u00d718=1152\", \"23\u00d769=1587\"),\n  @(\"53\u00d744=2332\", \"50\u00d745=2250\"),\n  @(\"27\u00d783=2241\", \"56\u00d782=4592\"),\n  @(\"28\u00d757=1596\", \"30\u00d740=1200\"),\n  @(\"23\u00d725=575\", \"54\u00d723=1242\"),\n  @(\"44\u00d778=3432\", \"34\u00d740=1360\"),\n  @(\"45\u00d798=4410\", \"17\u00d771=1207\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
